{"js": "const replacements = [\n  [\"828\u00f73=\", \"315\u00f73=\"],\n  [\"269\u00f79=\", \"497\u00f78=\"],\n  [\"699\u00f77=\", \"118\u00f78=\"],\n  [\"195\u00f74=\", \"267\u00f77=\"],\n  [\"263\u00f72=\", \"750\u00f75=\"],\n  [\"992\u00f73=\", \"500\u00f76=\"],\n  [\"129\u00f76=\", \"191\u00f78=\"],\n  [\"325\u00f78=\", \"417\u00f78=\"],\n  [\"250\u00f77=\", \"369\u00f75=\"],\n  [\"272\u00f74=\", \"795\u00f77=\"],\n  [\"436\u00f74=\", \"608\u00f73=\"],\n  [\"526\u00f79=\", \"208\u00f77=\"],\n  [\"527\u00f76=\", \"840\u00f78=\"],\n  [\"598\u00f77=\", \"172\u00f78=\"],\n  [\"683\u00f74=\", \"520\u00f73=\"],\n  [\"739\u00f77=\", \"198\u00f77=\"],\n  [\"654\u00f77=\", \"130\u00f75=\"],\n  [\"894\u00f73=\", \"590\u00f79=\"],\n  [\"735\u00f79=\", \"927\u00f77=\"],\n  [\"107\u00f79=\", \"957\u00f79=\"],\n  [\"646\u00f74=\", \"991\u00f75=\"],\n  [\"498\u00f75=\", \"849\u00f78=\"],\n  [\"221\u00f78=\", \"944\u00f77=\"],\n  [\"930\u00f72=\", \"779\u00f79=\"],\n  [\"881\u00f73=\", \"210\u00f72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"828\u00f73=\", \"315\u00f73=\"),\n  @(\"269\u00f79=\", \"497\u00f78=\"),\n  @(\"699\u00f77=\", \"118\u00f78=\"),\n  @(\"195\u00f74=\", \"267\u00f77=\"),\n  @(\"263\u00f72=\", \"750\u00f75=\"),\n  @(\"992\u00f73=\", \"500\u00f76=\"),\n  @(\"129\u00f76=\", \"191\u00f78=\"),\n  @(\"325\u00f78=\", \"417\u00f78=\"),\n  @(\"250\u00f77=\", \"369\u00f75=\"),\n  @(\"272\u00f74=\", \"795\u00f77=\"),\n  @(\"436\u00f74=\", \"608\u00f73=\"),\n  @(\"526\u00f79=\", \"208\u00f77=\"),\n  @(\"527\u00f76=\", \"840\u00f78=\"),\n  @(\"598\u00f77=\", \"172\u00f78=\"),\n  @(\"683\u00f74=\", \"520\u00f73=\"),\n  @(\"739\u00f77=\", \"198\u00f77=\"),\n  @(\"654\u00f77=\", \"130\u00f75=\"),\n  @(\"894\u00f73=\", \"590\u00f79=\"),\n  @(\"735\u00f79=\", \"927\u00f77=\"),\n  @(\"107\u00f79=\", \"957\u00f79=\"),\n  @(\"646\u00f74=\", \"991\u00f75=\"),\n  @(\"498\u00f75=\", \"849\u00f78=\"),\n  @(\"221\u00f78=\", \"944\u00f77=\"),\n  @(\"930\u00f72=\", \"779\u00f79=\"),\n  @(\"881\u00f73=\", \"210\u00f72=\"),\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
